# Fix the font size and alignment in the "Google Shared Drive to Google
# Shared Drive" advanced-features table: every paragraph/run in that
# table should end up at 10.5pt (w:sz/w:szCs = 21 half-points) for both
# the Latin and complex-script font sizes, including the paragraph mark
# run properties (the <w:rPr> nested inside <w:pPr>).
#
# Word's object model exposes this as Font.Size (points, drives w:sz)
# and Font.SizeBi (points, drives w:szCs). Setting them across the
# whole table Range touches every run in every cell as well as each
# paragraph's mark formatting, which is exactly the shape of the diff.

$d = $word.ActiveDocument

$table = $d.Tables(2)

$table.Range.Font.Size = 10.5
$table.Range.Font.SizeBi = 10.5
